# Add multisite beta-diversity indices (sp.sim, sp.sne, sp.sor, gen.sim,
# gen.sne, gen.sor, sf.sim, sf.sne, sf.sor) as new columns BI:BQ on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column headers, written into row 1 starting at column 61 (BI)
$headers = @("sp.sim","sp.sne","sp.sor","gen.sim","gen.sne","gen.sor","sf.sim","sf.sne","sf.sor")
$firstCol = 61

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $firstCol + $i).Value = $headers[$i]
}

# New data values per data row (rows without a value keep no entry here,
# matching rows 7, 12, 19 and 20 which have no computed result upstream)
$rowData = @{
    2 = @(0.63912164114417802, 0.26650024188438498, 0.90562188302856295, 0.38286713286713298, 0.50384660269127501, 0.88671373555840804, 0.041237113402061903, 0.76818911610613505, 0.80942622950819698);
    3 = @(0.62353775257001098, 0.24633113604658, 0.86986888861659095, 0.47217806041335503, 0.355085207016406, 0.82726326742976097, 0, 0.693965517241379, 0.693965517241379);
    4 = @(0.62104514752932805, 0.23929809890064699, 0.86034324642997495, 0.493365500603136, 0.28423577024674501, 0.77760127084988095, 0.028846153846153799, 0.60979427549194998, 0.63864042933810405);
    5 = @(0.63092269326683303, 0.19124612619847201, 0.82216881946530496, 0.38768115942029002, 0.35070893345896698, 0.73839009287925705, 0.06, 0.60784452296819802, 0.66784452296819796);
    6 = @(0.82718204488778102, 0.096820696716057697, 0.924002741603838, 0.66917859834212501, 0.21747257645777199, 0.88665117479989697, 0.0204081632653061, 0.57538131041890395, 0.59578947368421098);
    8 = @(0.65354330708661401, 0.17390767330554299, 0.82745098039215703, 0.50704225352112697, 0.219520246478873, 0.7265625, 0, 0.53333333333333299, 0.53333333333333299);
    9 = @(0.66305525460454995, 0.11839239121414299, 0.78144764581869297, 0.125, 0.32327586206896602, 0.44827586206896602, 0, 0.29310344827586199, 0.29310344827586199);
    10 = @(0.36216216216216202, 0.294315275974665, 0.65647743813682702, 0.2, 0.37088122605363999, 0.57088122605364, 0, 0.44144144144144098, 0.44144144144144098);
    11 = @(0.71348314606741603, 0.16144021383448501, 0.87492335990190095, 0.58628841607564997, 0.23940879507972801, 0.82569721115537897, 0.074074074074074098, 0.57381324986958804, 0.647887323943662);
    13 = @(0.57538461538461505, 0.22001938758054401, 0.79540400296515901, 0.271356783919598, 0.42273604308462098, 0.69409282700421904, 0, 0.57068062827225097, 0.57068062827225097);
    14 = @(0.60634081902245696, 0.22548762341546599, 0.83182844243792298, 0.35666666666666702, 0.43377850162866399, 0.79044516829533096, 0, 0.71317829457364301, 0.71317829457364301);
    15 = @(0.28571428571428598, 0.32212885154061599, 0.60784313725490202, 0.108108108108108, 0.35963382737576299, 0.467741935483871, 0, 0.21052631578947401, 0.21052631578947401);
    16 = @(0.38501291989664099, 0.229249640881317, 0.61426256077795804, 0.12258064516129, 0.29372836771424599, 0.41630901287553601, 0, 0.095238095238095205, 0.095238095238095205);
    17 = @(0.49270482603815902, 0.333348666361302, 0.82605349239946102, 0.253676470588235, 0.52674537906022301, 0.78042184964845895, 0, 0.64261168384879697, 0.64261168384879697);
    18 = @(0.45299145299145299, 0.25320594563134702, 0.70619739862279995, 0.16666666666666699, 0.38050314465408802, 0.54716981132075504, 0, 0.36666666666666697, 0.36666666666666697);
    21 = @(0.316114109483423, 0.41304619586008801, 0.72916030534351195, 0.057142857142857099, 0.581411359724613, 0.63855421686747005, 0, 0.471014492753623, 0.471014492753623)
}

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($r, $firstCol + $i).Value = $vals[$i]
    }
}

# Move the selection in the bottom-right (frozen) pane onto the newly
# entered block of data, mirroring the author's on-screen selection.
$ws.Range("BI8:BQ11").Select()
